# Swahili (Kenya) translation pass for
# "Facilitators guidelines - Game of Life.docx"
#
# Replaces the English table-header / label strings with their Swahili
# translations. Each Find/Replace is scoped to the whole document body,
# uses exact case matching and whole-word matching so we only ever touch
# the intended (unique) run in the document.

$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    $range = $d.Content
    $found = $range.Find.Execute(
        $oldText,   # FindText
        $true,      # MatchCase
        $true,      # MatchWholeWord
        $false,     # MatchWildcards
        $false,     # MatchSoundsLike
        $false,     # MatchAllWordForms
        $true,      # Forward
        1,          # Wrap -> wdFindContinue
        $false,     # Format
        $newText,   # ReplaceWith
        2           # Replace -> wdReplaceAll
    )
    if (-not $found) {
        Write-Output "WARNING: text not found -> '$oldText'"
    }
}

Replace-ExactText "Video Title" "Kichwa cha Video"
Replace-ExactText "Topic" "Mada"
Replace-ExactText "Aim(s)" "Malengo"
Replace-ExactText "Length" "Urefu"
Replace-ExactText "Camp Location" "Mahali pa Kambi"
Replace-ExactText "Facilitators" "Wawezeshaji"
Replace-ExactText "N. of students" "N. ya wanafunzi"
Replace-ExactText "Date" "Tarehe"
Replace-ExactText "Resources" "Rasilimali"
Replace-ExactText "needed" "inahitajika"
Replace-ExactText "Preparations" "Maandalizi"
Replace-ExactText "Video time" "Muda wa video"
Replace-ExactText "What facilitator does" "Mwezeshaji anafanya nini"
Replace-ExactText "What learners do" "Wanachofanya wanafunzi"
Replace-ExactText "General VMC Video Introduction" "Utangulizi Mkuu wa Video ya VMC"
Replace-ExactText "Video Introduction" "Utangulizi wa Video"
Replace-ExactText "Assist the process, provoke thoughts" "Kusaidia mchakato, kuchochea mawazo"
Replace-ExactText "Solution" "Suluhisho"

Write-Output "done"
